# Updated policies and graphs
# Set the "Industries" policy flag (column H) to 0 for rows 25 through 66.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("H25:H66").Value = 0
